# Update the "Promocionales Concentrados" report:
#  - refresh the cut-off date in the title
#  - fill row 3 with the latest promo-week record
#  - drop the now-unused trailing blank rows (9-12)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Promocionales Concentrados")

# --- helper: write a value as genuine text (t="s"), without leaving any
#     quote-prefix / number-format style behind on the cell -------------
function Set-TextValue {
    param($cell, [string]$text)

    $escaped = $text -replace '"', '""'
    $cell.Formula = '="' + $escaped + '"'
    $cell.Copy()
    $cell.PasteSpecial(-4163)  # xlPasteValues
    $excel.CutCopyMode = 0
}

# 1) Title cell (C1) — new cut-off date ----------------------------------
$ws.Range("C1").Value = "Informe Semana Promocional Consentrados a Corte: 08 - junio - 2022"

# 2) Row 3 — new data record --------------------------------------------
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "SO6117"
$ws.Range("C3").Value = "FEAL0000000102"
Set-TextValue $ws.Range("D3") "910000081"
$ws.Range("E3").Value = "CHUNKY GATOS POLLO 8 KG"
$ws.Range("F3").Value = 84340
$ws.Range("G3").Value = 84340
Set-TextValue $ws.Range("H3") "0"
Set-TextValue $ws.Range("I3") "1"
$ws.Range("J3").Value = $false
Set-TextValue $ws.Range("K3") "1"
$ws.Range("L3").Value = 910
Set-TextValue $ws.Range("N3") "2022-05-19"

# 3) Drop the now-empty trailing rows 9-12 -> dimension shrinks to N8 ----
$ws.Range("A9:N12").EntireRow.Delete()
